{"js": "// Replace the date line and all 25 \"a\u00d7b=c\" table-cell answers with the\n// updated values from the target revision. Every original text run in\n// the document is unique, so a direct search + full-text replace for\n// each pair is unambiguous and keeps each run's original formatting\n// (font, size, etc.) intact.\nconst replacements = [\n  [\"2024-04-04 Thursday\", \"2024-04-05 Friday\"],\n  [\"719\u00d78=5752\", \"933\u00d74=3732\"],\n  [\"128\u00d75=640\", \"738\u00d79=6642\"],\n  [\"631\u00d76=3786\", \"161\u00d74=644\"],\n  [\"447\u00d79=4023\", \"518\u00d76=3108\"],\n  [\"802\u00d76=4812\", \"567\u00d79=5103\"],\n  [\"252\u00d76=1512\", \"841\u00d79=7569\"],\n  [\"482\u00d78=3856\", \"283\u00d77=1981\"],\n  [\"944\u00d77=6608\", \"359\u00d75=1795\"],\n  [\"900\u00d79=8100\", \"256\u00d74=1024\"],\n  [\"324\u00d76=1944\", \"171\u00d79=1539\"],\n  [\"401\u00d76=2406\", \"997\u00d75=4985\"],\n  [\"755\u00d73=2265\", \"600\u00d72=1200\"],\n  [\"169\u00d77=1183\", \"906\u00d78=7248\"],\n  [\"584\u00d78=4672\", \"281\u00d76=1686\"],\n  [\"750\u00d77=5250\", \"685\u00d79=6165\"],\n  [\"299\u00d73=897\", \"378\u00d72=756\"],\n  [\"837\u00d76=5022\", \"942\u00d74=3768\"],\n  [\"828\u00d75=4140\", \"955\u00d76=5730\"],\n  [\"373\u00d79=3357\", \"350\u00d79=3150\"],\n  [\"797\u00d78=6376\", \"913\u00d75=4565\"],\n  [\"463\u00d75=2315\", \"939\u00d74=3756\"],\n  [\"541\u00d74=2164\", \"224\u00d78=1792\"],\n  [\"209\u00d76=1254\", \"998\u00d75=4990\"],\n  [\"670\u00d78=5360\", \"366\u00d77=2562\"],\n  [\"904\u00d74=3616\", \"304\u00d76=1824\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all 25 \"a\u00d7b=c\" table-cell answers with the\n# updated values from the target revision. Every original text string\n# is unique in the document, so Find/Replace (wdReplaceAll ==> 2, but\n# one hit each so a single execute suffices) for each pair swaps only\n# the intended run's text and leaves its formatting untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-04-04 Thursday\", \"2024-04-05 Friday\"),\n  @(\"719\u00d78=5752\", \"933\u00d74=3732\"),\n  @(\"128\u00d75=640\", \"738\u00d79=6642\"),\n  @(\"631\u00d76=3786\", \"161\u00d74=644\"),\n  @(\"447\u00d79=4023\", \"518\u00d76=3108\"),\n  @(\"802\u00d76=4812\", \"567\u00d79=5103\"),\n  @(\"252\u00d76=1512\", \"841\u00d79=7569\"),\n  @(\"482\u00d78=3856\", \"283\u00d77=1981\"),\n  @(\"944\u00d77=6608\", \"359\u00d75=1795\"),\n  @(\"900\u00d79=8100\", \"256\u00d74=1024\"),\n  @(\"324\u00d76=1944\", \"171\u00d79=1539\"),\n  @(\"401\u00d76=2406\", \"997\u00d75=4985\"),\n  @(\"755\u00d73=2265\", \"600\u00d72=1200\"),\n  @(\"169\u00d77=1183\", \"906\u00d78=7248\"),\n  @(\"584\u00d78=4672\", \"281\u00d76=1686\"),\n  @(\"750\u00d77=5250\", \"685\u00d79=6165\"),\n  @(\"299\u00d73=897\", \"378\u00d72=756\"),\n  @(\"837\u00d76=5022\", \"942\u00d74=3768\"),\n  @(\"828\u00d75=4140\", \"955\u00d76=5730\"),\n  @(\"373\u00d79=3357\", \"350\u00d79=3150\"),\n  @(\"797\u00d78=6376\", \"913\u00d75=4565\"),\n  @(\"463\u00d75=2315\", \"939\u00d74=3756\"),\n  @(\"541\u00d74=2164\", \"224\u00d78=1792\"),\n  @(\"209\u00d76=1254\", \"998\u00d75=4990\"),\n  @(\"670\u00d78=5360\", \"366\u00d77=2562\"),\n  @(\"904\u00d74=3616\", \"304\u00d76=1824\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
